$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows (628-633) for a new price-reporting week (2022-02-03, serial 44595)
$ws.Range("A628:R633").Insert()

$ws.Range("A628").Value = 2
$ws.Range("B628").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C628").Value = "Coquimbo"
$ws.Range("D628").Value = 44595
$ws.Range("E628").Value = 4
$ws.Range("F628").Value = 100112002
$ws.Range("G628").Value = "Pimiento"
$ws.Range("H628").Value = "Cuatro cascos rojo"
$ws.Range("I628").Value = "Primera"
$ws.Range("J628").Value = 300
$ws.Range("K628").Value = 11000
$ws.Range("L628").Value = 12000
$ws.Range("M628").Value = 11500
$ws.Range("N628").Value = "$/caja 18 kilos"
$ws.Range("O628").Value = "Provincia de Limarí"
$ws.Range("P628").Value = 639
$ws.Range("Q628").Value = 18
$ws.Range("R628").Value = "Hortaliza"

$ws.Range("A629").Value = 2
$ws.Range("B629").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C629").Value = "Coquimbo"
$ws.Range("D629").Value = 44595
$ws.Range("E629").Value = 4
$ws.Range("F629").Value = 100112002
$ws.Range("G629").Value = "Pimiento"
$ws.Range("H629").Value = "Cuatro cascos rojo"
$ws.Range("I629").Value = "Segunda"
$ws.Range("J629").Value = 300
$ws.Range("K629").Value = 8000
$ws.Range("L629").Value = 9000
$ws.Range("M629").Value = 8500
$ws.Range("N629").Value = "$/caja 18 kilos"
$ws.Range("O629").Value = "Provincia de Limarí"
$ws.Range("P629").Value = 472
$ws.Range("Q629").Value = 18
$ws.Range("R629").Value = "Hortaliza"

$ws.Range("A630").Value = 2
$ws.Range("B630").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C630").Value = "Coquimbo"
$ws.Range("D630").Value = 44595
$ws.Range("E630").Value = 4
$ws.Range("F630").Value = 100112002
$ws.Range("G630").Value = "Pimiento"
$ws.Range("H630").Value = "Cuatro cascos rojo"
$ws.Range("I630").Value = "Tercera"
$ws.Range("J630").Value = 300
$ws.Range("K630").Value = 4000
$ws.Range("L630").Value = 5000
$ws.Range("M630").Value = 4500
$ws.Range("N630").Value = "$/caja 18 kilos"
$ws.Range("O630").Value = "Provincia de Limarí"
$ws.Range("P630").Value = 250
$ws.Range("Q630").Value = 18
$ws.Range("R630").Value = "Hortaliza"

$ws.Range("A631").Value = 2
$ws.Range("B631").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C631").Value = "Coquimbo"
$ws.Range("D631").Value = 44595
$ws.Range("E631").Value = 4
$ws.Range("F631").Value = 100112002
$ws.Range("G631").Value = "Pimiento"
$ws.Range("H631").Value = "Cuatro cascos verde"
$ws.Range("I631").Value = "Primera"
$ws.Range("J631").Value = 400
$ws.Range("K631").Value = 7500
$ws.Range("L631").Value = 8000
$ws.Range("M631").Value = 7750
$ws.Range("N631").Value = "$/caja 18 kilos"
$ws.Range("O631").Value = "Provincia de Limarí"
$ws.Range("P631").Value = 431
$ws.Range("Q631").Value = 18
$ws.Range("R631").Value = "Hortaliza"

$ws.Range("A632").Value = 2
$ws.Range("B632").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C632").Value = "Coquimbo"
$ws.Range("D632").Value = 44595
$ws.Range("E632").Value = 4
$ws.Range("F632").Value = 100112002
$ws.Range("G632").Value = "Pimiento"
$ws.Range("H632").Value = "Cuatro cascos verde"
$ws.Range("I632").Value = "Segunda"
$ws.Range("J632").Value = 360
$ws.Range("K632").Value = 5500
$ws.Range("L632").Value = 6000
$ws.Range("M632").Value = 5750
$ws.Range("N632").Value = "$/caja 18 kilos"
$ws.Range("O632").Value = "Provincia de Limarí"
$ws.Range("P632").Value = 319
$ws.Range("Q632").Value = 18
$ws.Range("R632").Value = "Hortaliza"

$ws.Range("A633").Value = 2
$ws.Range("B633").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C633").Value = "Coquimbo"
$ws.Range("D633").Value = 44595
$ws.Range("E633").Value = 4
$ws.Range("F633").Value = 100112002
$ws.Range("G633").Value = "Pimiento"
$ws.Range("H633").Value = "Cuatro cascos verde"
$ws.Range("I633").Value = "Tercera"
$ws.Range("J633").Value = 200
$ws.Range("K633").Value = 3500
$ws.Range("L633").Value = 4000
$ws.Range("M633").Value = 3750
$ws.Range("N633").Value = "$/caja 18 kilos"
$ws.Range("O633").Value = "Provincia de Limarí"
$ws.Range("P633").Value = 208
$ws.Range("Q633").Value = 18
$ws.Range("R633").Value = "Hortaliza"

